$d = $word.ActiveDocument

# Update the date heading
$d.Paragraphs.Item(1).Range.Text = "2023-08-21 Monday"

# Update the division problems in the table (row, col) -> new value
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "34÷4="  # was 29÷3=
$t.Cell(1,2).Range.Text = "26÷7="  # was 81÷6=
$t.Cell(1,3).Range.Text = "15÷6="  # was 15÷5=
$t.Cell(1,4).Range.Text = "15÷4="  # was 80÷9=
$t.Cell(1,5).Range.Text = "53÷7="  # was 21÷9=
$t.Cell(5,1).Range.Text = "30÷5="  # was 46÷2=
$t.Cell(5,2).Range.Text = "11÷4="  # was 27÷4=
$t.Cell(5,3).Range.Text = "14÷5="  # was 18÷7=
$t.Cell(5,4).Range.Text = "96÷5="  # was 82÷9=
$t.Cell(5,5).Range.Text = "25÷9="  # was 69÷5=
$t.Cell(9,1).Range.Text = "49÷4="  # was 40÷9=
$t.Cell(9,2).Range.Text = "71÷8="  # was 51÷6=
$t.Cell(9,3).Range.Text = "18÷7="  # was 63÷4=
$t.Cell(9,4).Range.Text = "18÷3="  # was 21÷9=
$t.Cell(9,5).Range.Text = "36÷2="  # was 61÷2=
$t.Cell(13,1).Range.Text = "72÷6="  # was 57÷5=
$t.Cell(13,2).Range.Text = "11÷5="  # was 88÷4=
$t.Cell(13,3).Range.Text = "88÷2="  # was 91÷7=
$t.Cell(13,4).Range.Text = "61÷2="  # was 65÷6=
$t.Cell(13,5).Range.Text = "44÷9="  # was 48÷7=
$t.Cell(17,1).Range.Text = "75÷6="  # was 16÷9=
$t.Cell(17,2).Range.Text = "46÷9="  # was 42÷2=
$t.Cell(17,3).Range.Text = "44÷9="  # was 27÷8=
$t.Cell(17,4).Range.Text = "30÷4="  # was 50÷3=
$t.Cell(17,5).Range.Text = "45÷4="  # was 75÷6=
